# Remove the first paragraph of the document body, which contains the
# "I personally examined the patient separately ..." attestation text.
# The diff removes this paragraph (including its paragraph mark)
# entirely, so the document now starts with the "OBJECTIVE:" paragraph.

$d = $word.ActiveDocument

$target = $d.Paragraphs(1)
$target.Range.Delete()
